$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.659958333333333
$ws.Range("H2").Value = 10.979875
$ws.Range("I2").Value = 0.4781132044744068
$ws.Range("J2").Value = 0.4781132044744067
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.046397
$ws.Range("N2").Value = 0.139191
$ws.Range("O2").Value = 0.01970724914037141
$ws.Range("P2").Value = 0.01970724914037141
$ws.Range("Q2").Value = 0.1698110867916667
$ws.Range("R2").Value = 1.528299781125
$ws.Range("S2").Value = 0.009422296037878473
$ws.Range("T2").Value = 0.009422296037878471
$ws.Range("G3").Value = 3.659958333333333
$ws.Range("H3").Value = 10.979875
$ws.Range("I3").Value = 0.4781132044744068
$ws.Range("J3").Value = 0.4781132044744067
$ws.Range("O3").Value = 0.6598912010221247
$ws.Range("P3").Value = 0.6598912010221247
$ws.Range("Q3").Value = 5.686072227111111
$ws.Range("R3").Value = 51.174650044
$ws.Range("S3").Value = 0.315502696725153
$ws.Range("T3").Value = 0.3155026967251529
$ws.Range("G4").Value = 3.659958333333333
$ws.Range("H4").Value = 10.979875
$ws.Range("I4").Value = 0.4781132044744068
$ws.Range("J4").Value = 0.4781132044744067
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.754325
$ws.Range("N4").Value = 2.262975
$ws.Range("O4").Value = 0.3204015498375038
$ws.Range("P4").Value = 0.3204015498375038
$ws.Range("Q4").Value = 2.760798069791667
$ws.Range("R4").Value = 24.847182628125
$ws.Range("S4").Value = 0.1531882117113753
$ws.Range("T4").Value = 0.1531882117113753
$ws.Range("I5").Value = 0.3193330932870009
$ws.Range("J5").Value = 0.3193330932870008
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.046397
$ws.Range("N5").Value = 0.139191
$ws.Range("O5").Value = 0.01970724914037141
$ws.Range("P5").Value = 0.01970724914037141
$ws.Range("Q5").Value = 0.113417280912
$ws.Range("R5").Value = 1.020755528208
$ws.Range("S5").Value = 0.006293176828172391
$ws.Range("T5").Value = 0.00629317682817239
$ws.Range("I6").Value = 0.3193330932870009
$ws.Range("J6").Value = 0.3193330932870008
$ws.Range("O6").Value = 0.6598912010221247
$ws.Range("P6").Value = 0.6598912010221247
$ws.Range("S6").Value = 0.2107250984552692
$ws.Range("T6").Value = 0.2107250984552692
$ws.Range("I7").Value = 0.3193330932870009
$ws.Range("J7").Value = 0.3193330932870008
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.754325
$ws.Range("N7").Value = 2.262975
$ws.Range("O7").Value = 0.3204015498375038
$ws.Range("P7").Value = 0.3204015498375038
$ws.Range("Q7").Value = 1.8439444452
$ws.Range("R7").Value = 16.5955000068
$ws.Range("S7").Value = 0.1023148180035593
$ws.Range("T7").Value = 0.1023148180035592
$ws.Range("G8").Value = 1.388093333333333
$ws.Range("H8").Value = 4.16428
$ws.Range("I8").Value = 0.1813315046964271
$ws.Range("J8").Value = 0.1813315046964271
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.046397
$ws.Range("N8").Value = 0.139191
$ws.Range("O8").Value = 0.01970724914037141
$ws.Range("P8").Value = 0.01970724914037141
$ws.Range("Q8").Value = 0.06440336638666666
$ws.Range("R8").Value = 0.5796302974800001
$ws.Range("S8").Value = 0.003573545140050917
$ws.Range("T8").Value = 0.003573545140050917
$ws.Range("G9").Value = 1.388093333333333
$ws.Range("H9").Value = 4.16428
$ws.Range("I9").Value = 0.1813315046964271
$ws.Range("J9").Value = 0.1813315046964271
$ws.Range("O9").Value = 0.6598912010221247
$ws.Range("P9").Value = 0.6598912010221247
$ws.Range("Q9").Value = 2.156526996337778
$ws.Range("R9").Value = 19.40874296704
$ws.Range("S9").Value = 0.1196590644172743
$ws.Range("T9").Value = 0.1196590644172743
$ws.Range("G10").Value = 1.388093333333333
$ws.Range("H10").Value = 4.16428
$ws.Range("I10").Value = 0.1813315046964271
$ws.Range("J10").Value = 0.1813315046964271
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.754325
$ws.Range("N10").Value = 2.262975
$ws.Range("O10").Value = 0.3204015498375038
$ws.Range("P10").Value = 0.3204015498375038
$ws.Range("Q10").Value = 1.047073503666667
$ws.Range("R10").Value = 9.423661532999999
$ws.Range("S10").Value = 0.05809889513910185
$ws.Range("T10").Value = 0.05809889513910184
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.162456
$ws.Range("H11").Value = 0.487368
$ws.Range("I11").Value = 0.02122219754216535
$ws.Range("J11").Value = 0.02122219754216534
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.046397
$ws.Range("N11").Value = 0.139191
$ws.Range("O11").Value = 0.01970724914037141
$ws.Range("P11").Value = 0.01970724914037141
$ws.Range("Q11").Value = 0.007537471032000001
$ws.Range("R11").Value = 0.06783723928800001
$ws.Range("S11").Value = 0.0004182311342696302
$ws.Range("T11").Value = 0.0004182311342696302
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.162456
$ws.Range("H12").Value = 0.487368
$ws.Range("I12").Value = 0.02122219754216535
$ws.Range("J12").Value = 0.02122219754216534
$ws.Range("O12").Value = 0.6598912010221247
$ws.Range("P12").Value = 0.6598912010221247
$ws.Range("Q12").Value = 0.252389908736
$ws.Range("R12").Value = 2.271509178624
$ws.Range("S12").Value = 0.01400434142442827
$ws.Range("T12").Value = 0.01400434142442827
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.162456
$ws.Range("H13").Value = 0.487368
$ws.Range("I13").Value = 0.02122219754216535
$ws.Range("J13").Value = 0.02122219754216534
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.754325
$ws.Range("N13").Value = 2.262975
$ws.Range("O13").Value = 0.3204015498375038
$ws.Range("P13").Value = 0.3204015498375038
$ws.Range("Q13").Value = 0.1225446222
$ws.Range("R13").Value = 1.1029015998
$ws.Range("S13").Value = 0.006799624983467441
$ws.Range("T13").Value = 0.006799624983467439
